$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + "28.426.09"
$ws.Range("E2").Value = "  +0.03%  "

$ws.Range("D3").Value = "'" + "1.553.83"
$ws.Range("E3").Value = "  -1.91%  "

$ws.Range("D5").Value = "'" + "210.53"
$ws.Range("E5").Value = "  -1.35%  "

$ws.Range("E6").Value = "  -1.79%  "

$ws.Range("D8").Value = "'" + "24.12"
$ws.Range("E8").Value = "  +0.84%  "

$ws.Range("E9").Value = "  -2.03%  "

$ws.Range("D10").Value = "'" + "0.0583"
$ws.Range("E10").Value = "  -1.20%  "

$ws.Range("D11").Value = "'" + "0.0891"
$ws.Range("E11").Value = "  -0.40%  "

$ws.Range("D12").Value = "'" + "1.775.70"
$ws.Range("E12").Value = "  -2.03%  "

$ws.Range("D13").Value = "'" + "1.557.43"
$ws.Range("E13").Value = "  -1.51%  "

$ws.Range("D14").Value = "'" + "28.449.46"
$ws.Range("E14").Value = "  -0.03%  "

$ws.Range("D15").Value = "'" + "3.63"
$ws.Range("E15").Value = "  -2.20%  "

$ws.Range("E16").Value = "  -1.59%  "

$ws.Range("E17").Value = "  -1.76%  "

$ws.Range("D18").Value = "'" + "228.99"
$ws.Range("E18").Value = "  -0.45%  "

$ws.Range("E19").Value = "  -1.55%  "

$ws.Range("D20").Value = "'" + "0.0" + [char]0x2083 + "0672"
$ws.Range("E20").Value = "  -2.47%  "

$ws.Range("E21").Value = "  -0.15%  "

$ws.Range("D22").Value = "'" + "3.88"
$ws.Range("E22").Value = "  -0.73%  "

$ws.Range("E23").Value = "  -2.61%  "

$ws.Range("E24").Value = "  -2.27%  "

$ws.Range("D25").Value = "'" + "150.95"
$ws.Range("E25").Value = "  -0.47%  "

$ws.Range("E27").Value = "  -1.38%  "

$ws.Range("E28").Value = "  -0.21%  "

$ws.Range("D29").Value = "'" + "6.23"
$ws.Range("E29").Value = "  -3.16%  "

$ws.Range("E30").Value = "  -3.16%  "

$ws.Range("E31").Value = "  -4.48%  "

$ws.Range("E32").Value = "  -1.68%  "

$ws.Range("D33").Value = "'" + "1.384.01"
$ws.Range("E33").Value = "  -0.74%  "

$ws.Range("E34").Value = "  -3.02%  "

$ws.Range("E35").Value = "  -2.89%  "

$ws.Range("D36").Value = "'" + "1.48"
$ws.Range("E36").Value = "  -2.16%  "

$ws.Range("E37").Value = "  -2.81%  "

$ws.Range("E38").Value = "  -3.07%  "

$ws.Range("E39").Value = "  -2.46%  "

$ws.Range("D40").Value = "'" + "1.93"
$ws.Range("E40").Value = "  +2.49%  "

$ws.Range("D41").Value = "'" + "0.510"
$ws.Range("E41").Value = "  -2.25%  "

$ws.Range("E42").Value = "  -0.18%  "

$ws.Range("D43").Value = "'" + "0.772"
$ws.Range("E43").Value = "  -2.29%  "

$ws.Range("D44").Value = "'" + "0.0461"
$ws.Range("E44").Value = "  +0.82%  "

$ws.Range("E45").Value = "  -1.54%  "

$ws.Range("E46").Value = "  -1.85%  "

$ws.Range("D47").Value = "'" + "1.688.14"
$ws.Range("E47").Value = "  -1.97%  "

$ws.Range("D48").Value = "'" + "0.874"
$ws.Range("E48").Value = "  -8.86%  "

$ws.Range("D49").Value = "'" + "85.05"
$ws.Range("E49").Value = "  -1.82%  "

$ws.Range("D50").Value = "'" + "42.93"
$ws.Range("E50").Value = "  +7.54%  "

$ws.Range("E51").Value = "  -1.35%  "
